$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of column L (rows 3-12) into column M, then set the new values.
$ws.Range("L3:L12").Copy()
$ws.Range("M3:M12").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("M4").Value = 2023
$ws.Range("M5").Value = 311.65582791395695
$ws.Range("M7").Value = 119.55977988994496
$ws.Range("M8").Value = 192.09604802401199
$ws.Range("M10").Value = 78.539269634817401
$ws.Range("M11").Value = 60.030015007503756
$ws.Range("M12").Value = 26.013006503251628

$ws.Range("A1").Select()
